$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Continuing of Chapter 3  (Documentation)"
$ws.Range("D6").Value = "8pm"
$ws.Range("E6").Value = "1pm"

$ws.Range("B5:C5").Copy()
$ws.Range("B6:C6").PasteSpecial(-4122)

$ws.Range("B6").Value2 = 43755
$ws.Range("C6").Value2 = 43755

$ws.Range("A6").Select()
